$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update row 2 - person record
$ws.Range("A2").Value = "Monica"
$ws.Range("B2").Value = "Francisca"
$ws.Range("C2").Value = "Gastanbide"
$ws.Range("D2").Value = 2025240810
$ws.Range("E2").Value = "MonicaF1"

# Update row 3 - person record
$ws.Range("A3").Value = "Charles"
$ws.Range("B3").Value = "Antonio"
$ws.Range("C3").Value = "Lecrec Montez"
$ws.Range("D3").Value = 2025240811
$ws.Range("E3").Value = "CharlesAL112"

# Update the selected cell to F3
$ws.Range("F3").Select()

$wb.Save()
